# Remove the stray "<lb/>" run that immediately follows the "</env>" run
# in the paragraph "<env>a la pluye</env><lb/></ab>". The search string
# "</env><lb/>" is unique across the whole document, and since the match
# starts exactly at a run boundary, replacing it with "</env>" preserves
# the formatting (Courier New, blue, sz 18) of the existing "</env>" run
# while the trailing "<lb/>" run (Courier New, gray, sz 18) is dropped
# entirely — matching the target diff.
$d = $word.ActiveDocument
$d.Content.Find.Execute("</env><lb/>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "</env>", 2)
